$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect the new "through" date
$ws.Name = "Through 2022-04-14"

# Update the row label for April
$ws.Range("A5").Value = "April (through 04-14)"

# Update April row (row 5) values for the columns that changed
$ws.Range("C5").Value = 11
$ws.Range("D5").Value = 26
$ws.Range("E5").Value = 23
$ws.Range("F5").Value = 23
$ws.Range("H5").Value = 49
$ws.Range("I5").Value = 59

# Update Total row (row 6) values for the columns that changed
$ws.Range("C6").Value = 139
$ws.Range("D6").Value = 215
$ws.Range("E6").Value = 220
$ws.Range("F6").Value = 133
$ws.Range("H6").Value = 472
$ws.Range("I6").Value = 493
